$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: ID "20" (reuse the existing shared string), Value 42
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4163)
$ws.Range("B9").Value = 42

# Row 10: ID "12" (new shared string), Value 12
# Build the text "12" via a formula returning a string, then paste its
# value so the destination cell ends up as a genuine text cell (not a
# formula) without ever touching NumberFormat (which would otherwise
# leave a stray style behind).
$ws.Range("Z1").Formula = "=""12"""
$ws.Range("Z1").Copy()
$ws.Range("A10").PasteSpecial(-4163)
$ws.Range("B10").Value = 12
$ws.Range("Z1").Clear()
